$wb = $excel.ActiveWorkbook

# Rename worksheets: replace "TRE-R" with "TRE_R" in the three "Include from TRE-..." sheets
$wb.Worksheets.Item("Include from TRE-R67-TypeStru").Name = "Include from TRE_R67-TypeStru"
$wb.Worksheets.Item("Include from TRE-R04-TypeSavo").Name = "Include from TRE_R04-TypeSavo"
$wb.Worksheets.Item("Include from TRE-R288-TypePro").Name = "Include from TRE_R288-TypePro"

# Update the Date value on the Metadata sheet (row 8: A8="Date", B8=timestamp)
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B8").Value = "2024-04-10T09:05:12+00:00"
